# state specific fix (#46)
# The "CDO Override reason" column header had inconsistent capitalization;
# fix it to "CDO Override Reason".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AF1").Value = "CDO Override Reason"

# Leave the fixed cell selected, matching the author's last edit position.
[void]$ws.Range("AF1").Select()
